$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "6th july content" values into column G for rows 3-6 (previously empty)
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0

# Extend the SUM formulas in column H to include the new column G values
$ws.Range("H3").Formula = "=SUM(B3:G3)"
$ws.Range("H4").Formula = "=SUM(B4:G4)"
$ws.Range("H5").Formula = "=SUM(B5:G5)"
$ws.Range("H6").Formula = "=SUM(B6:G6)"

# Move the active selection to H6, matching the saved view state
$ws.Range("H6").Select()
